$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header labels in B1 and C1 (shared strings)
$ws.Range("B1").Value = "AVERAGE_TIME_PER_ISSUES_WO_FT"
$ws.Range("C1").Value = "AVERAGE_TIME_PER_ISSUES_WITH_FT"

# Update the active cell / selection to B7
$ws.Range("B7").Select()
